$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(55, 1).Value = "I have diarrhea"
$ws.Cells.Item(55, 2).Value = "下痢です。|げりです。"
$ws.Cells.Item(56, 1).Value = "I am constipated."
$ws.Cells.Item(56, 2).Value = "便秘です。|べんぴです。"
$ws.Cells.Item(57, 1).Value = "I have my period."
$ws.Cells.Item(57, 2).Value = "生理です。|せいりです。"
$ws.Cells.Item(58, 1).Value = "I have hay fever."
$ws.Cells.Item(58, 2).Value = "花粉症です。|かふんしょうです。"
$ws.Cells.Item(59, 1).Value = "I have an allergy to..."
$ws.Cells.Item(59, 2).Value = "（～に）アレルギーがあります。"
$ws.Cells.Item(60, 1).Value = "I have a bad tooth."
$ws.Cells.Item(60, 2).Value = "虫歯があります。|むしばがあります。"
$ws.Cells.Item(61, 1).Value = "I sneeze."
$ws.Cells.Item(61, 2).Value = "くしゃみが出ます。|くしゃみがでます。"
$ws.Cells.Item(62, 1).Value = "I have a runny nose."
$ws.Cells.Item(62, 2).Value = "鼻水が出ます。|はなみずがでます。"
$ws.Cells.Item(63, 1).Value = "My back itches."
$ws.Cells.Item(63, 2).Value = "背中がかゆいです。|せなかがかゆいです。"
$ws.Cells.Item(64, 1).Value = "I have rashes."
$ws.Cells.Item(64, 2).Value = "発疹があります。|はっしんがあります。"
$ws.Cells.Item(65, 1).Value = "I feel dizzy."
$ws.Cells.Item(65, 2).Value = "めまいがします。"
$ws.Cells.Item(66, 1).Value = "I threw up."
$ws.Cells.Item(66, 2).Value = "吐きました。|はきました。"
$ws.Cells.Item(67, 1).Value = "I am not feeling well."
$ws.Cells.Item(67, 2).Value = "気分が悪いです。|きぶんがわるいです。"
$ws.Cells.Item(68, 1).Value = "I burned myself."
$ws.Cells.Item(68, 2).Value = "やけどをしました。"
$ws.Cells.Item(69, 1).Value = "I broke my leg."
$ws.Cells.Item(69, 2).Value = "足を骨折しました。|あしをこっせつしました。"
$ws.Cells.Item(70, 1).Value = "I hurt myself."
$ws.Cells.Item(70, 2).Value = "けがをしました。"
$ws.Cells.Item(71, 1).Value = "physician"
$ws.Cells.Item(71, 2).Value = "内科|ないか"
$ws.Cells.Item(72, 1).Value = "dermatologist"
$ws.Cells.Item(72, 2).Value = "皮膚科|ひふか"
$ws.Cells.Item(73, 1).Value = "surgeon"
$ws.Cells.Item(73, 2).Value = "外科|げか"
$ws.Cells.Item(74, 1).Value = "obstetrician and gynecologist"
$ws.Cells.Item(74, 2).Value = "産婦人科|さんふじんか"
$ws.Cells.Item(75, 1).Value = "orthopedic surgeon"
$ws.Cells.Item(75, 2).Value = "整形外科|せいけいげか"
$ws.Cells.Item(76, 1).Value = "ophthalmologist"
$ws.Cells.Item(76, 2).Value = "眼科|がんか"
$ws.Cells.Item(77, 1).Value = "dentist"
$ws.Cells.Item(77, 2).Value = "歯科|しか"
$ws.Cells.Item(78, 1).Value = "otorhinolaryngologist; ENT doctor"
$ws.Cells.Item(78, 2).Value = "耳鼻科|じびか"
$ws.Cells.Item(79, 1).Value = "antibiotic"
$ws.Cells.Item(79, 2).Value = "抗生物質|こうせいぶっしつ"
$ws.Cells.Item(80, 1).Value = "X-ray"
$ws.Cells.Item(80, 2).Value = "レントゲン"
$ws.Cells.Item(81, 1).Value = "operation"
$ws.Cells.Item(81, 2).Value = "手術|しゅじゅつ"
$ws.Cells.Item(82, 1).Value = "injection"
$ws.Cells.Item(82, 2).Value = "注射|ちゅうしゃ"
$ws.Cells.Item(83, 1).Value = "thermometer"
$ws.Cells.Item(83, 2).Value = "体温計|たいおんけい"
$ws.Cells.Item(84, 1).Value = "intravenous feeding"
$ws.Cells.Item(84, 2).Value = "点滴|てんてき"
